# Homepage / calendar UI modification:
# - Add a new "venue_name" column (D) with a placeholder "-" value for every match row.
# - Recolor body/header font to an explicit black (was theme-based black).
# - Bump every row's height slightly to fit the new column's content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count   # 25 -> header + 24 data rows

# 1) New header cell for the venue_name column.
$ws.Range("D1").Value2 = "venue_name"

# 2) Fill the rest of column D with the "-" placeholder (rows 2..last).
$dataRange = $ws.Range("D2:D" + $lastRow)
$dataRange.Value2 = "-"

# 3) Make the font color explicit black (rgb) rather than theme-derived black,
#    across the whole populated area (header + data, columns A:E).
$usedRange = $ws.Range("A1:E" + $lastRow)
$usedRange.Font.Color = 0

# 4) Row heights grow slightly now that every row carries the extra column.
$newHeights = @{
  1=19.5; 2=20.25; 3=19.5; 4=19.5; 5=19.5; 6=20.25; 7=19.5; 8=19.5; 9=19.5; 10=20.25;
  11=20.25; 12=19.5; 13=19.5; 14=19.5; 15=19.5; 16=19.5; 17=19.5; 18=19.5; 19=19.5; 20=19.5;
  21=19.5; 22=19.5; 23=19.5; 24=19.5; 25=20.25
}
foreach ($r in $newHeights.Keys) {
  $ws.Rows.Item($r).RowHeight = $newHeights[$r]
}

# 5) Columns A-C get a left-aligned default (previously "general"),
#    matching the rest of the table's alignment.
$ws.Columns("A:C").HorizontalAlignment = -4131
